$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 186.66667
$ws.Range("I6").Value = 191.25
$ws.Range("K6").Value = 573.75
$ws.Range("M6").Value = -461.75
$ws.Range("H28").Value = 1587.1111
$ws.Range("I28").Value = 1879.1666
$ws.Range("M28").Value = -1394.1666
$ws.Range("K28").Value = 1879.1666
$ws.Range("J51").Value = 25499.334
$ws.Range("H51").Value = 25499.334
$ws.Range("N51").Value = -26467.334
$ws.Range("L51").Value = 25499.334
$ws.Range("L62").Value = 2900
$ws.Range("J62").Value = 2900
$ws.Range("H62").Value = 2900
$ws.Range("N62").Value = -4148
$ws.Range("H64").Value = 8166.923
$ws.Range("J64").Value = 9263.223
$ws.Range("N64").Value = -9759.223
$ws.Range("L64").Value = 9263.223
$ws.Range("L65").Value = 14500
$ws.Range("N65").Value = -20740
$ws.Range("H65").Value = 2900
$ws.Range("J65").Value = 2900
$ws.Range("L67").Value = 9263.223
$ws.Range("N67").Value = -10979.223
$ws.Range("J67").Value = 9263.223
$ws.Range("H67").Value = 8166.923
$ws.Range("K96").Value = 9078
$ws.Range("H96").Value = 3513
$ws.Range("M96").Value = -7705
$ws.Range("I96").Value = 3026
$ws.Range("I107").Value = 1045
$ws.Range("H107").Value = 1375.4
$ws.Range("L107").Value = 1595.6666
$ws.Range("M107").Value = 875
$ws.Range("N107").Value = -5435.6666
$ws.Range("K107").Value = 1045
$ws.Range("J107").Value = 1595.6666
$ws.Range("I111").Value = 4696.8335
$ws.Range("M111").Value = -11023.5005
$ws.Range("K111").Value = 14090.5005
$ws.Range("H111").Value = 7454.2856
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("L130").Value = 120000
$ws.Range("N130").Value = -130040
$ws.Range("H130").Value = 120000
$ws.Range("J130").Value = 120000
$ws.Range("M135").Value = -21531.5625
$ws.Range("H135").Value = 7698.8423
$ws.Range("K135").Value = 24066.5625
$ws.Range("I135").Value = 2674.0625
$ws.Range("I137").Value = 1999.6923
$ws.Range("K137").Value = 5999.0769
$ws.Range("H137").Value = 3524.15
$ws.Range("M137").Value = -3449.0769
$ws.Range("J138").Value = 5321.1177
$ws.Range("I138").Value = 1570.0344
$ws.Range("K138").Value = 4710.1032
$ws.Range("N138").Value = -26243.3531
$ws.Range("L138").Value = 15963.3531
$ws.Range("M138").Value = 429.8968000000004
$ws.Range("H138").Value = 3961.35

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L2").Value = 4447
$ws.Range("N2").Value = -4673
$ws.Range("H2").Value = 3931.9333
$ws.Range("J2").Value = 4447
$ws.Range("K21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("I21").Value = 0
$ws.Range("K25").Value = 783
$ws.Range("H25").Value = 783
$ws.Range("I25").Value = 783
$ws.Range("M25").Value = -381
$ws.Range("I32").Value = 813.0833
$ws.Range("M32").Value = -526.0833
$ws.Range("K32").Value = 813.0833
$ws.Range("H32").Value = 7943.3335
$ws.Range("M45").Value = -19446.818
$ws.Range("H45").Value = 13902.471
$ws.Range("N45").Value = -3800.6667
$ws.Range("L45").Value = 3046.6667
$ws.Range("K45").Value = 19823.818
$ws.Range("J45").Value = 3046.6667
$ws.Range("I45").Value = 19823.818
$ws.Range("M61").Value = -4213.7188
$ws.Range("J61").Value = 4734.5
$ws.Range("I61").Value = 4425.7188
$ws.Range("H61").Value = 4460.028
$ws.Range("L61").Value = 4734.5
$ws.Range("K61").Value = 4425.7188
$ws.Range("N61").Value = -5158.5
$ws.Range("L74").Value = 2797
$ws.Range("N74").Value = -4545
$ws.Range("K74").Value = 0
$ws.Range("H74").Value = 2797
$ws.Range("M74").ClearContents()
$ws.Range("J74").Value = 2797
$ws.Range("I74").Value = 0
$ws.Range("J77").Value = 2797
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 13985
$ws.Range("N77").Value = -22721
$ws.Range("H77").Value = 2797
$ws.Range("M77").ClearContents()
$ws.Range("I77").Value = 0
$ws.Range("N94").Value = -141790.38
$ws.Range("H94").Value = 139988.38
$ws.Range("J94").Value = 139988.38
$ws.Range("L94").Value = 139988.38
$ws.Range("M97").Value = -632
$ws.Range("H97").Value = 1433.8379
$ws.Range("I97").Value = 1128
$ws.Range("K97").Value = 1128
$ws.Range("J105").Value = 73333.336
$ws.Range("H105").Value = 73333.336
$ws.Range("N105").Value = -80321.336
$ws.Range("L105").Value = 73333.336
$ws.Range("J110").Value = 881.2
$ws.Range("N110").Value = -4971.2
$ws.Range("H110").Value = 1262.8125
$ws.Range("K110").Value = 1436.2727
$ws.Range("M110").Value = 608.7273
$ws.Range("L110").Value = 881.2
$ws.Range("I110").Value = 1436.2727
$ws.Range("L116").Value = 4447
$ws.Range("N116").Value = -9035
$ws.Range("J116").Value = 4447
$ws.Range("H116").Value = 3931.9333
$ws.Range("K122").Value = 7958.571599999999
$ws.Range("J122").Value = 1800
$ws.Range("N122").Value = -10300
$ws.Range("M122").Value = -5508.571599999999
$ws.Range("I122").Value = 2652.8572
$ws.Range("L122").Value = 5400
$ws.Range("H122").Value = 2546.25
$ws.Range("N128").Value = -83293
$ws.Range("H128").Value = 73333
$ws.Range("L128").Value = 73333
$ws.Range("J128").Value = 73333
$ws.Range("N134").Value = -110140
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 100000
$ws.Range("H134").Value = 100000
$ws.Range("J136").Value = 4734.5
$ws.Range("K136").Value = 13277.1564
$ws.Range("H136").Value = 4460.028
$ws.Range("M136").Value = -10727.1564
$ws.Range("L136").Value = 14203.5
$ws.Range("I136").Value = 4425.7188
$ws.Range("N136").Value = -19303.5
$ws.Range("J138").Value = 115000
$ws.Range("N138").Value = -125280
$ws.Range("L138").Value = 115000
$ws.Range("H138").Value = 115000
$ws.Range("H139").Value = 75000
$ws.Range("N139").Value = -85280
$ws.Range("J139").Value = 75000
$ws.Range("L139").Value = 75000
$ws.Range("N141").Value = -63359.668
$ws.Range("H141").Value = 52999.668
$ws.Range("J141").Value = 52999.668
$ws.Range("L141").Value = 52999.668

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J3").Value = 4447
$ws.Range("N3").Value = -4675
$ws.Range("L3").Value = 4447
$ws.Range("H3").Value = 3931.9333
$ws.Range("H19").Value = 4044.4443
$ws.Range("K19").Value = 4044.4443
$ws.Range("M19").Value = -3871.4443
$ws.Range("I19").Value = 4044.4443
$ws.Range("H22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L80").Value = 2121.3333
$ws.Range("J80").Value = 2121.3333
$ws.Range("K80").Value = 1586.5
$ws.Range("N80").Value = -4117.3333
$ws.Range("I80").Value = 1586.5
$ws.Range("H80").Value = 1907.4
$ws.Range("M80").Value = -588.5
$ws.Range("M83").Value = -2940.5
$ws.Range("I83").Value = 1586.5
$ws.Range("K83").Value = 7932.5
$ws.Range("H83").Value = 1907.4
$ws.Range("L83").Value = 10606.6665
$ws.Range("J83").Value = 2121.3333
$ws.Range("N83").Value = -20590.6665
$ws.Range("I94").Value = 656.2143
$ws.Range("M94").Value = -205.2143
$ws.Range("K94").Value = 656.2143
$ws.Range("H94").Value = 681.05884
$ws.Range("H99").Value = 3055.75
$ws.Range("I99").Value = 2842.05
$ws.Range("K99").Value = 2842.05
$ws.Range("M99").Value = -1344.05
$ws.Range("M105").Value = -7475
$ws.Range("K105").Value = 9222
$ws.Range("I105").Value = 9222
$ws.Range("H105").Value = 9222
$ws.Range("I107").Value = 1312.6562
$ws.Range("H107").Value = 1452.3096
$ws.Range("M107").Value = 607.3438000000001
$ws.Range("K107").Value = 1312.6562
$ws.Range("H125").Value = 83748.25
$ws.Range("L125").Value = 83748.25
$ws.Range("N125").Value = -93588.25
$ws.Range("J125").Value = 83748.25
$ws.Range("N126").Value = -129879.5
$ws.Range("J126").Value = 119999.5
$ws.Range("H126").Value = 119999.5
$ws.Range("L126").Value = 119999.5
$ws.Range("M134").Value = -5647.399800000001
$ws.Range("I134").Value = 2727.4666
$ws.Range("H134").Value = 3772.2778
$ws.Range("K134").Value = 8182.399800000001
$ws.Range("N137").Value = -75590
$ws.Range("H137").Value = 65390
$ws.Range("J137").Value = 65390
$ws.Range("L137").Value = 65390

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I16").Value = 15692.9
$ws.Range("H16").Value = 20462.424
$ws.Range("L16").Value = 27800.154
$ws.Range("K16").Value = 15692.9
$ws.Range("N16").Value = -28374.154
$ws.Range("J16").Value = 27800.154
$ws.Range("M16").Value = -15405.9
$ws.Range("J31").Value = 4987
$ws.Range("I31").Value = 1978.6
$ws.Range("K31").Value = 1978.6
$ws.Range("N31").Value = -5577
$ws.Range("H31").Value = 2922.4119
$ws.Range("M31").Value = -1683.6
$ws.Range("L31").Value = 4987
$ws.Range("H34").Value = 2922.4119
$ws.Range("J34").Value = 4987
$ws.Range("I34").Value = 1978.6
$ws.Range("N34").Value = -5391
$ws.Range("L34").Value = 4987
$ws.Range("M34").Value = -1776.6
$ws.Range("K34").Value = 1978.6
$ws.Range("K39").Value = 6520.3335
$ws.Range("N39").Value = -40781
$ws.Range("H39").Value = 14890
$ws.Range("L39").Value = 39999
$ws.Range("I39").Value = 6520.3335
$ws.Range("M39").Value = -6129.3335
$ws.Range("J39").Value = 39999
$ws.Range("N49").Value = -40363
$ws.Range("I49").Value = 6520.3335
$ws.Range("H49").Value = 14890
$ws.Range("L49").Value = 39999
$ws.Range("K49").Value = 6520.3335
$ws.Range("M49").Value = -6338.3335
$ws.Range("J49").Value = 39999
$ws.Range("L58").Value = 4474.75
$ws.Range("M58").Value = -7239.6665
$ws.Range("I58").Value = 7442.6665
$ws.Range("H58").Value = 6817.8423
$ws.Range("K58").Value = 7442.6665
$ws.Range("J58").Value = 4474.75
$ws.Range("N58").Value = -4880.75
$ws.Range("N59").Value = -64789.75
$ws.Range("H59").Value = 46571.285
$ws.Range("J59").Value = 62499.75
$ws.Range("L59").Value = 62499.75
$ws.Range("I62").Value = 4731.857
$ws.Range("K62").Value = 4731.857
$ws.Range("H62").Value = 5400.9414
$ws.Range("M62").Value = -4107.857
$ws.Range("M65").Value = -20539.285
$ws.Range("H65").Value = 5400.9414
$ws.Range("K65").Value = 23659.285
$ws.Range("I65").Value = 4731.857
$ws.Range("K86").Value = 20839136
$ws.Range("H86").Value = 17550242
$ws.Range("M86").Value = -20838013
$ws.Range("I86").Value = 20839136
$ws.Range("H89").Value = 17550242
$ws.Range("I89").Value = 20839136
$ws.Range("M89").Value = -104190064
$ws.Range("K89").Value = 104195680
$ws.Range("N98").Value = -84450
$ws.Range("L98").Value = 79958
$ws.Range("J98").Value = 79958
$ws.Range("H98").Value = 79958
$ws.Range("N113").Value = -32140.154
$ws.Range("M113").Value = -13522.9
$ws.Range("I113").Value = 15692.9
$ws.Range("H113").Value = 20462.424
$ws.Range("J113").Value = 27800.154
$ws.Range("K113").Value = 15692.9
$ws.Range("L113").Value = 27800.154
$ws.Range("K122").Value = 6288.1875
$ws.Range("M122").Value = -3838.1875
$ws.Range("I122").Value = 2096.0625
$ws.Range("H122").Value = 3056.2173
$ws.Range("M132").Value = -6941.643199999999
$ws.Range("I132").Value = 3157.2144
$ws.Range("K132").Value = 9471.643199999999
$ws.Range("H132").Value = 4300
$ws.Range("M134").Value = -3591
$ws.Range("I134").Value = 2042
$ws.Range("H134").Value = 2772.0527
$ws.Range("K134").Value = 6126
$ws.Range("J136").Value = 4474.75
$ws.Range("K136").Value = 22327.9995
$ws.Range("H136").Value = 6817.8423
$ws.Range("M136").Value = -19777.9995
$ws.Range("L136").Value = 13424.25
$ws.Range("I136").Value = 7442.6665
$ws.Range("N136").Value = -18524.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K3").Value = 10236
$ws.Range("I3").Value = 3412
$ws.Range("M3").Value = -10124
$ws.Range("H3").Value = 3412
$ws.Range("L40").Value = 590
$ws.Range("H40").Value = 138
$ws.Range("J40").Value = 147.5
$ws.Range("N40").Value = -728
$ws.Range("H50").Value = 1302.5625
$ws.Range("N50").Value = -6025.727000000001
$ws.Range("L50").Value = 5063.727000000001
$ws.Range("J50").Value = 1687.909
$ws.Range("N53").Value = -6025.727000000001
$ws.Range("L53").Value = 5063.727000000001
$ws.Range("J53").Value = 1687.909
$ws.Range("H53").Value = 1302.5625
$ws.Range("H101").Value = 20029
$ws.Range("L101").Value = 60087
$ws.Range("N101").Value = -64955
$ws.Range("J101").Value = 20029
$ws.Range("N113").Value = -13939.7501
$ws.Range("H113").Value = 2356.3809
$ws.Range("J113").Value = 3199.9167
$ws.Range("L113").Value = 9599.750100000001
$ws.Range("M134").Value = 3715.73685
$ws.Range("I134").Value = 451.42105
$ws.Range("H134").Value = 451.42105
$ws.Range("K134").Value = 1354.26315
$ws.Range("I138").Value = 3810
$ws.Range("K138").Value = 11430
$ws.Range("M138").Value = -6290
$ws.Range("H138").Value = 10285.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7204000.5
$ws.Range("M11").Value = -3147719.2
$ws.Range("I11").Value = 3147858.2
$ws.Range("K11").Value = 3147858.2
$ws.Range("L70").Value = 4900
$ws.Range("N70").Value = -5440
$ws.Range("H70").Value = 4902.3335
$ws.Range("J70").Value = 4900
$ws.Range("L73").Value = 4900
$ws.Range("H73").Value = 4902.3335
$ws.Range("N73").Value = -6772
$ws.Range("J73").Value = 4900
$ws.Range("L92").Value = 8501
$ws.Range("H92").Value = 8501
$ws.Range("N92").Value = -12245
$ws.Range("J92").Value = 8501
$ws.Range("M97").Value = 139.6
$ws.Range("H97").Value = 332.92856
$ws.Range("I97").Value = 356.4
$ws.Range("K97").Value = 356.4
$ws.Range("N103").ClearContents()
$ws.Range("L103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("H103").Value = 0
$ws.Range("I107").Value = 231.72728
$ws.Range("H107").Value = 304.07693
$ws.Range("L107").Value = 357.13333
$ws.Range("M107").Value = 1688.27272
$ws.Range("N107").Value = -4197.13333
$ws.Range("K107").Value = 231.72728
$ws.Range("J107").Value = 357.13333
$ws.Range("K122").Value = 9645.714
$ws.Range("J122").Value = 6626.067
$ws.Range("N122").Value = -24778.201
$ws.Range("M122").Value = -7195.714
$ws.Range("I122").Value = 3215.238
$ws.Range("L122").Value = 19878.201
$ws.Range("H122").Value = 4636.4165
$ws.Range("N133").Value = -90010
$ws.Range("H133").Value = 79890
$ws.Range("L133").Value = 79890
$ws.Range("J133").Value = 79890
$ws.Range("H135").Value = 84740
$ws.Range("L135").Value = 84740
$ws.Range("N135").Value = -94880
$ws.Range("J135").Value = 84740
$ws.Range("H139").Value = 99999.5
$ws.Range("N139").Value = -110279.5
$ws.Range("J139").Value = 99999.5
$ws.Range("L139").Value = 99999.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K7").Value = 12550.091
$ws.Range("I7").Value = 12550.091
$ws.Range("H7").Value = 13058
$ws.Range("M7").Value = -12438.091
$ws.Range("L9").Value = 14199.8
$ws.Range("J9").Value = 14199.8
$ws.Range("N9").Value = -14647.8
$ws.Range("M9").Value = -4177.4
$ws.Range("H9").Value = 9300.6
$ws.Range("K9").Value = 4401.4
$ws.Range("I9").Value = 4401.4
$ws.Range("I16").Value = 2791.1428
$ws.Range("H16").Value = 2860.6924
$ws.Range("L16").Value = 2941.8333
$ws.Range("K16").Value = 2791.1428
$ws.Range("N16").Value = -3281.8333
$ws.Range("J16").Value = 2941.8333
$ws.Range("M16").Value = -2621.1428
$ws.Range("J18").Value = 94998.664
$ws.Range("N18").Value = -95342.664
$ws.Range("M18").Value = -90828
$ws.Range("L18").Value = 94998.664
$ws.Range("H18").Value = 92999.336
$ws.Range("K18").Value = 91000
$ws.Range("I18").Value = 91000
$ws.Range("H22").Value = 2309.75
$ws.Range("M22").Value = -2008.75
$ws.Range("I22").Value = 2303.75
$ws.Range("K22").Value = 2303.75
$ws.Range("M27").Value = -2196.75
$ws.Range("K27").Value = 2303.75
$ws.Range("I27").Value = 2303.75
$ws.Range("H27").Value = 2309.75
$ws.Range("M40").Value = -2858.158
$ws.Range("L40").Value = 14433.875
$ws.Range("K40").Value = 2994.158
$ws.Range("H40").Value = 6383.7036
$ws.Range("I40").Value = 2994.158
$ws.Range("J40").Value = 14433.875
$ws.Range("N40").Value = -14705.875
$ws.Range("I46").Value = 6778
$ws.Range("M46").Value = -6590
$ws.Range("K46").Value = 6778
$ws.Range("H46").Value = 5963.1665
$ws.Range("L55").Value = 350.57144
$ws.Range("I55").Value = 175.83333
$ws.Range("H55").Value = 269.92307
$ws.Range("K55").Value = 175.83333
$ws.Range("N55").Value = -696.5714399999999
$ws.Range("M55").Value = -2.833329999999989
$ws.Range("J55").Value = 350.57144
$ws.Range("N68").Value = -4663.6667
$ws.Range("H68").Value = 2913.2856
$ws.Range("L68").Value = 3165.6667
$ws.Range("J68").Value = 3165.6667
$ws.Range("L71").Value = 15828.3335
$ws.Range("H71").Value = 2913.2856
$ws.Range("N71").Value = -23316.3335
$ws.Range("J71").Value = 3165.6667
$ws.Range("H100").Value = 2580.6
$ws.Range("K100").Value = 1949.5
$ws.Range("M100").Value = -1408.5
$ws.Range("I100").Value = 1949.5
$ws.Range("K122").Value = 10950.6
$ws.Range("J122").Value = 3998.75
$ws.Range("N122").Value = -16896.25
$ws.Range("M122").Value = -8500.599999999999
$ws.Range("I122").Value = 3650.2
$ws.Range("L122").Value = 11996.25
$ws.Range("H122").Value = 3723.5789
$ws.Range("K126").Value = 37650.273
$ws.Range("I126").Value = 12550.091
$ws.Range("H126").Value = 13058
$ws.Range("M126").Value = -35180.273
$ws.Range("J127").Value = 49951.855
$ws.Range("L127").Value = 49951.855
$ws.Range("N127").Value = -59871.855
$ws.Range("H127").Value = 49957.875
$ws.Range("N128").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("M132").Value = -86935.427
$ws.Range("I132").Value = 29821.809
$ws.Range("K132").Value = 89465.427
$ws.Range("H132").Value = 22097.19
$ws.Range("K136").Value = 13332
$ws.Range("H136").Value = 4721.75
$ws.Range("M136").Value = -10782
$ws.Range("I136").Value = 4444
$ws.Range("N137").Value = -90198.5
$ws.Range("H137").Value = 79998.5
$ws.Range("J137").Value = 79998.5
$ws.Range("L137").Value = 79998.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I62").Value = 4108.6665
$ws.Range("L62").Value = 11288.25
$ws.Range("J62").Value = 11288.25
$ws.Range("K62").Value = 4108.6665
$ws.Range("H62").Value = 6317.769
$ws.Range("M62").Value = -3484.6665
$ws.Range("N62").Value = -12536.25
$ws.Range("L65").Value = 56441.25
$ws.Range("M65").Value = -17423.3325
$ws.Range("N65").Value = -62681.25
$ws.Range("H65").Value = 6317.769
$ws.Range("J65").Value = 11288.25
$ws.Range("K65").Value = 20543.3325
$ws.Range("I65").Value = 4108.6665
$ws.Range("K96").Value = 940.6
$ws.Range("H96").Value = 931.8333
$ws.Range("M96").Value = 432.4
$ws.Range("I96").Value = 940.6
$ws.Range("H100").Value = 1002.625
$ws.Range("K100").Value = 2063.1428
$ws.Range("M100").Value = -1522.1428
$ws.Range("I100").Value = 1031.5714
$ws.Range("H107").Value = 2062.926
$ws.Range("L107").Value = 9721.799999999999
$ws.Range("N107").Value = -13561.8
$ws.Range("J107").Value = 3240.6
$ws.Range("N113").Value = -9890
$ws.Range("M113").Value = -2522.9
$ws.Range("I113").Value = 1564.3
$ws.Range("H113").Value = 1590.2727
$ws.Range("J113").Value = 1850
$ws.Range("K113").Value = 4692.9
$ws.Range("L113").Value = 5550
$ws.Range("J119").Value = 29022.5
$ws.Range("L119").Value = 29022.5
$ws.Range("N119").Value = -38698.5
$ws.Range("H119").Value = 1223218
$ws.Range("K122").Value = 14853.75
$ws.Range("J122").Value = 7762.375
$ws.Range("N122").Value = -28187.125
$ws.Range("M122").Value = -12403.75
$ws.Range("I122").Value = 4951.25
$ws.Range("L122").Value = 23287.125
$ws.Range("H122").Value = 6075.7
$ws.Range("N126").Value = -25224.8
$ws.Range("K126").Value = 15264.834
$ws.Range("I126").Value = 5088.278
$ws.Range("J126").Value = 6761.6
$ws.Range("H126").Value = 5452.0435
$ws.Range("L126").Value = 20284.8
$ws.Range("M126").Value = -12794.834
$ws.Range("M132").Value = -782.6921000000002
$ws.Range("N132").Value = -10456.4
$ws.Range("L132").Value = 5396.4
$ws.Range("I132").Value = 1104.2307
$ws.Range("K132").Value = 3312.6921
$ws.Range("H132").Value = 1297.1666
$ws.Range("J132").Value = 1798.8
$ws.Range("K136").Value = 5391
$ws.Range("H136").Value = 1797
$ws.Range("M136").Value = -2841
$ws.Range("I136").Value = 1797
$ws.Range("N137").Value = -73998.60000000001
$ws.Range("H137").Value = 63798.6
$ws.Range("J137").Value = 63798.6
$ws.Range("L137").Value = 63798.6
$ws.Range("J138").Value = 123499.5
$ws.Range("N138").Value = -133779.5
$ws.Range("L138").Value = 123499.5
$ws.Range("H138").Value = 123499.5
